{"js": "// Expand storefront theme catalog and add category filters.\n// Appends a new \"Update\" block (blank separator line, heading line,\n// header row, and five data rows) to the end of the document body,\n// matching the existing \"Update ...\" block formatting/style already\n// used throughout the document (Helvetica Light, 24 half-points).\n\nconst newParagraphs = [\n  \"\",\n  \"Update 2026-02-19 08:48 IST - Storefront System Theme Marketplace\",\n  \"Module name | developed | partial developed | need to develop\",\n  \"Customizable themes/templates | 10 category themes (free + paid) seeded, plan-eligibility apply flow active, merchant-side theme filter/search in Store Builder UI | Theme previews are image/card-based; runtime rendering still section-json centric | Full packaged theme runtime engine with deeper template inheritance and advanced live preview\",\n  \"Mobile responsive design | Store builder and theme catalog UI supports responsive grid behavior | Per-theme mobile QA matrix not automated | Dedicated performance budgets + Lighthouse gates per theme variant\",\n  \"Homepage builder (drag & drop sections) | Existing section-based builder + reorder + validation + versions already available | Full visual freeform canvas still partial | Advanced WYSIWYG constraints and reusable section presets per theme\",\n  \"Navigation / menus | Existing primary menu JSON editor is available | No visual nested menu designer yet | Full multi-level visual menu editor with visibility rules by customer type/device\",\n  \"Static pages (About, Contact, Policy) | CRUD + SEO fields already available | Rich block templates for policy/legal pages not yet | Template library + compliance presets per category\"\n];\n\nconst body = context.document.body;\n\n// Grab the formatting (font + size) off the very last existing\n// paragraph's run so the newly appended paragraphs match the rest of\n// the document (Helvetica Light, 24 half-points / 12pt).\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nconst lastRange = lastParagraph.getRange(\"Whole\");\nconst font = lastRange.font;\nfont.load(\"name,size\");\nawait context.sync();\n\nconst fontName = font.name;\nconst fontSize = font.size;\n\nfor (const text of newParagraphs) {\n  const p = body.insertParagraph(text, Word.InsertLocation.end);\n  p.font.name = fontName;\n  p.font.size = fontSize;\n  if (text === \"\") {\n    // Force a (empty) <w:t> run to materialize for the blank\n    // separator paragraph, matching the rest of the document where\n    // every paragraph carries an explicit (possibly empty) text run.\n    p.getRange().insertText(\"\", Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Expand storefront theme catalog and add category filters.\n# Appends a new \"Update\" block (blank separator line, heading line,\n# header row, and five data rows) to the end of the document body,\n# matching the existing \"Update ...\" block formatting/style already\n# used throughout the document (Helvetica Light, 24 half-points / 12pt).\n\n$d = $word.ActiveDocument\n\n$newParagraphs = @(\n  \"\",\n  \"Update 2026-02-19 08:48 IST - Storefront System Theme Marketplace\",\n  \"Module name | developed | partial developed | need to develop\",\n  \"Customizable themes/templates | 10 category themes (free + paid) seeded, plan-eligibility apply flow active, merchant-side theme filter/search in Store Builder UI | Theme previews are image/card-based; runtime rendering still section-json centric | Full packaged theme runtime engine with deeper template inheritance and advanced live preview\",\n  \"Mobile responsive design | Store builder and theme catalog UI supports responsive grid behavior | Per-theme mobile QA matrix not automated | Dedicated performance budgets + Lighthouse gates per theme variant\",\n  \"Homepage builder (drag & drop sections) | Existing section-based builder + reorder + validation + versions already available | Full visual freeform canvas still partial | Advanced WYSIWYG constraints and reusable section presets per theme\",\n  \"Navigation / menus | Existing primary menu JSON editor is available | No visual nested menu designer yet | Full multi-level visual menu editor with visibility rules by customer type/device\",\n  \"Static pages (About, Contact, Policy) | CRUD + SEO fields already available | Rich block templates for policy/legal pages not yet | Template library + compliance presets per category\"\n)\n\n# No explicit font/size assignment needed: InsertParagraphAfter /\n# InsertAfter inherit the run formatting (Helvetica Light, 24\n# half-points) already in effect at the end of the document, which\n# matches every other paragraph in the file.\nforeach ($t in $newParagraphs) {\n  $r = $d.Content\n  $r.Collapse(0)\n  $r.InsertParagraphAfter()\n  $r.Collapse(0)\n  $r.InsertAfter($t)\n}\n"}
